$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Ariful Islam: F7/G7 0 -> 40 (now fully paid 80 hrs), I7 "Short 80.00 hours" -> blank
$ws.Range("F7").Value = 40
$ws.Range("G7").Value = 40
$ws.Range("I7").Value = ""

# Row 11 - Pauline Nguyen: F11/G11 0 -> 40 (now fully paid 80 hrs), I11 "Short 80.00 hours" -> blank
$ws.Range("F11").Value = 40
$ws.Range("G11").Value = 40
$ws.Range("I11").Value = ""

# Row 12 - Edward Obi: G12 0 -> 35 (F12 stays 0), I12 "Short 45.00 hours" -> "Short 10.00 hours"
$ws.Range("G12").Value = 35
$ws.Range("I12").Value = "Short 10.00 hours"
